# Auto-generated edit script applying the diff to before.xlsx
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item('展览')
$ws1.Range("F2").Value = 1941
$ws1.Range("F4").Value = 101
$ws1.Range("F5").Value = 430
$ws1.Range("F6").Value = 1846
$ws1.Range("F7").Value = 864
$ws1.Range("F8").Value = 1313
$ws1.Range("F9").Value = 565
$ws1.Range("F10").Value = 128
$ws1.Range("F11").Value = 2797
$ws1.Range("F12").Value = 370
$ws1.Range("F14").Value = 1099
$ws1.Range("F15").Value = 583
$ws1.Range("F16").Value = 322
$ws1.Range("F17").Value = 63
$ws1.Range("F18").Value = 1603
$ws1.Range("F19").Value = 332
$ws1.Range("F20").Value = 1244
$ws1.Range("F21").Value = 186
$ws1.Range("F22").Value = 596
$ws1.Range("F23").Value = 502
$ws1.Range("F25").Value = 1466
$ws1.Range("F26").Value = 1454
$ws1.Range("F27").Value = 1325
$ws1.Range("F28").Value = 269
$ws1.Range("F29").Value = 1277
$ws1.Range("F30").Value = 432
$ws1.Range("F31").Value = 146
$ws1.Range("F32").Value = 963
$ws1.Range("F33").Value = 24
$ws1.Range("F34").Value = 1841
$ws1.Range("F35").Value = 472
$ws1.Range("F36").Value = 39
$ws1.Range("F38").Value = 19
$ws1.Range("F39").Value = 2273
$ws1.Range("F40").Value = 145
$ws1.Range("F41").Value = 885
$ws1.Range("F42").Value = 2766
$ws1.Range("F43").Value = 11
$ws1.Range("F45").Value = 20

$ws2 = $wb.Worksheets.Item('演出')
$ws2.Range("F2").Value = 57
$ws2.Range("F3").Value = 61
$ws2.Range("F5").Value = 60
$ws2.Range("F7").Value = 28
$ws2.Range("F10").Value = 32
$ws2.Range("F12").Value = 363
$ws2.Range("F13").Value = 109583
$ws2.Range("F14").Value = 44
$ws2.Range("F17").Value = 66
$ws2.Range("F18").Value = 66
$ws2.Range("F20").Value = 287
$ws2.Range("F22").Value = 281
$ws2.Range("F24").Value = 80
$ws2.Range("F26").Value = 60
$ws2.Range("F27").Value = 60
$ws2.Range("F30").Value = 77
$ws2.Range("F34").Value = 17
$ws2.Range("F35").Value = 84

$ws3 = $wb.Worksheets.Item('本地生活')
$ws3.Range("F4").Value = 291
$ws3.Range("F6").Value = 4848
$ws3.Range("G6").Value = '已售罄'
$ws3.Range("F7").Value = 177
$ws3.Range("F9").Value = 663
$ws3.Range("F10").Value = 931
$ws3.Range("F11").Value = 541
$ws3.Range("F12").Value = 615
$ws3.Range("F13").Value = 1346
$ws3.Range("F14").Value = 384
$ws3.Range("F15").Value = 1229

$ws4 = $wb.Worksheets.Item('全部类型')
$ws4.Range("F2").Value = 1941
$ws4.Range("F3").Value = 291
$ws4.Range("B5").Value = "'2024-05-25"
$ws4.Range("C5").Value = '上海·「排球少年!!垃圾场决战 × animate cafe」'
$ws4.Range("D5").Value = '西藏北路198号大悦城北座8楼N809-1 animate cafe上海店'
$ws4.Range("E5").Value = '2024.05.25 00:00-07.02 23:59'
$ws4.Range("F5").Value = 4848
$ws4.Range("G5").Value = '已售罄'
$ws4.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=85283'
$ws4.Range("I5").Value = '//i0.hdslb.com/bfs/openplatform/202405/vy2vecK11715162037223.jpeg'
$ws4.Range("B6").Value = "'2024-06-01"
$ws4.Range("C6").Value = '上海·NIJISANJI EN 官方授权主题店'
$ws4.Range("D6").Value = '西藏北路166号（地铁8号线曲阜路下） 静安大悦城'
$ws4.Range("E6").Value = '2024.06.01 00:00-07.15 23:59'
$ws4.Range("F6").Value = 663
$ws4.Range("G6").Value = 30
$ws4.Range("H6").Value = 'https://show.bilibili.com/platform/detail.html?id=86310'
$ws4.Range("I6").Value = '//i2.hdslb.com/bfs/openplatform/202405/MhBVkfZ51716778059321.jpeg'
$ws4.Range("B7").Value = "'2024-06-07"
$ws4.Range("C7").Value = '上海·全职高手×HAPPY ZOO 全职高手十周年咖啡厅'
$ws4.Range("D7").Value = '南京东路340号百联zx创趣场四楼05号 HAPPY ZOO'
$ws4.Range("E7").Value = '2024.06.07 00:00-08.04 23:59'
$ws4.Range("F7").Value = 931
$ws4.Range("G7").Value = 10
$ws4.Range("H7").Value = 'https://show.bilibili.com/platform/detail.html?id=86871'
$ws4.Range("I7").Value = '//i0.hdslb.com/bfs/openplatform/202405/xw8aUE5u1715846379865.jpeg'
$ws4.Range("C8").Value = '上海· 怪兽8号 meets niko and … 集章之旅    '
$ws4.Range("D8").Value = '吴江路169号1层E127,E128 niko and ... (上海四季坊店)'
$ws4.Range("E8").Value = '2024.06.08 10:00-07.21 22:00'
$ws4.Range("F8").Value = 541
$ws4.Range("G8").Value = 48
$ws4.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=85758'
$ws4.Range("I8").Value = '//i0.hdslb.com/bfs/openplatform/202405/xw8aUE5u1715846379865.jpeg'
$ws4.Range("F9").Value = 615
$ws4.Range("F10").Value = 1346
$ws4.Range("F11").Value = 430
$ws4.Range("F12").Value = 1846
$ws4.Range("F13").Value = 864
$ws4.Range("F14").Value = 1313
$ws4.Range("F15").Value = 28
$ws4.Range("F16").Value = 565
$ws4.Range("F17").Value = 1229
$ws4.Range("F18").Value = 2797
$ws4.Range("F19").Value = 32
$ws4.Range("F20").Value = 370
$ws4.Range("F22").Value = 1099
$ws4.Range("F23").Value = 583
$ws4.Range("F24").Value = 322
$ws4.Range("F25").Value = 1603
$ws4.Range("F26").Value = 332
$ws4.Range("F27").Value = 363
$ws4.Range("F28").Value = 1244
$ws4.Range("F29").Value = 186
$ws4.Range("F30").Value = 596
$ws4.Range("F31").Value = 1466
$ws4.Range("F32").Value = 1454
$ws4.Range("F33").Value = 1325
$ws4.Range("F35").Value = 66
$ws4.Range("F36").Value = 1277
$ws4.Range("F37").Value = 432
$ws4.Range("F38").Value = 963
$ws4.Range("F40").Value = 1841
$ws4.Range("F41").Value = 60
$ws4.Range("F43").Value = 2273
$ws4.Range("F44").Value = 145
$ws4.Range("F45").Value = 885
$ws4.Range("F46").Value = 2766
